# Implement more stopwords for ngrams
#
# The "2017" sheet lists (n, ngram_text, ngram_count) rows computed from the
# headline "chinese lunar rover finds no evidence of american moon landings".
# Adding "no" and "of" to the ngram stopword list drops those two tokens, so
# the filtered sentence becomes
#   "chinese lunar rover finds evidence american moon landings"
# Every remaining n-gram shifts to a larger n (since two fewer stopwords are
# skipped over) and the 4-gram/5-gram rows that no longer have distinct
# n-grams (rows 32-39) disappear.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2017")

# Remove the trailing rows that no longer exist (n=4 had 7 rows, n=5 had 6 -
# both shrink once "no"/"of" stop being counted as ngram tokens).
$ws.Range("A32:C39").EntireRow.Delete()

$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = "rover"
$ws.Cells.Item(2,3).Value = 1
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "moon"
$ws.Cells.Item(3,3).Value = 1
$ws.Cells.Item(4,1).Value = 1
$ws.Cells.Item(4,2).Value = "lunar"
$ws.Cells.Item(4,3).Value = 1
$ws.Cells.Item(5,1).Value = 1
$ws.Cells.Item(5,2).Value = "landings"
$ws.Cells.Item(5,3).Value = 1
$ws.Cells.Item(6,1).Value = 1
$ws.Cells.Item(6,2).Value = "finds"
$ws.Cells.Item(6,3).Value = 1
$ws.Cells.Item(7,1).Value = 1
$ws.Cells.Item(7,2).Value = "evidence"
$ws.Cells.Item(7,3).Value = 1
$ws.Cells.Item(8,1).Value = 1
$ws.Cells.Item(8,2).Value = "chinese"
$ws.Cells.Item(8,3).Value = 1
$ws.Cells.Item(9,1).Value = 1
$ws.Cells.Item(9,2).Value = "american"
$ws.Cells.Item(9,3).Value = 1
$ws.Cells.Item(10,1).Value = 2
$ws.Cells.Item(10,2).Value = "rover finds"
$ws.Cells.Item(10,3).Value = 1
$ws.Cells.Item(11,1).Value = 2
$ws.Cells.Item(11,2).Value = "moon landings"
$ws.Cells.Item(11,3).Value = 1
$ws.Cells.Item(12,1).Value = 2
$ws.Cells.Item(12,2).Value = "lunar rover"
$ws.Cells.Item(12,3).Value = 1
$ws.Cells.Item(13,1).Value = 2
$ws.Cells.Item(13,2).Value = "finds evidence"
$ws.Cells.Item(13,3).Value = 1
$ws.Cells.Item(14,1).Value = 2
$ws.Cells.Item(14,2).Value = "evidence american"
$ws.Cells.Item(14,3).Value = 1
$ws.Cells.Item(15,1).Value = 2
$ws.Cells.Item(15,2).Value = "chinese lunar"
$ws.Cells.Item(15,3).Value = 1
$ws.Cells.Item(16,1).Value = 2
$ws.Cells.Item(16,2).Value = "american moon"
$ws.Cells.Item(16,3).Value = 1
$ws.Cells.Item(17,1).Value = 3
$ws.Cells.Item(17,2).Value = "rover finds evidence"
$ws.Cells.Item(17,3).Value = 1
$ws.Cells.Item(18,1).Value = 3
$ws.Cells.Item(18,2).Value = "lunar rover finds"
$ws.Cells.Item(18,3).Value = 1
$ws.Cells.Item(19,1).Value = 3
$ws.Cells.Item(19,2).Value = "finds evidence american"
$ws.Cells.Item(19,3).Value = 1
$ws.Cells.Item(20,1).Value = 3
$ws.Cells.Item(20,2).Value = "evidence american moon"
$ws.Cells.Item(20,3).Value = 1
$ws.Cells.Item(21,1).Value = 3
$ws.Cells.Item(21,2).Value = "chinese lunar rover"
$ws.Cells.Item(21,3).Value = 1
$ws.Cells.Item(22,1).Value = 3
$ws.Cells.Item(22,2).Value = "american moon landings"
$ws.Cells.Item(22,3).Value = 1
$ws.Cells.Item(23,1).Value = 4
$ws.Cells.Item(23,2).Value = "rover finds evidence american"
$ws.Cells.Item(23,3).Value = 1
$ws.Cells.Item(24,1).Value = 4
$ws.Cells.Item(24,2).Value = "lunar rover finds evidence"
$ws.Cells.Item(24,3).Value = 1
$ws.Cells.Item(25,1).Value = 4
$ws.Cells.Item(25,2).Value = "finds evidence american moon"
$ws.Cells.Item(25,3).Value = 1
$ws.Cells.Item(26,1).Value = 4
$ws.Cells.Item(26,2).Value = "evidence american moon landings"
$ws.Cells.Item(26,3).Value = 1
$ws.Cells.Item(27,1).Value = 4
$ws.Cells.Item(27,2).Value = "chinese lunar rover finds"
$ws.Cells.Item(27,3).Value = 1
$ws.Cells.Item(28,1).Value = 5
$ws.Cells.Item(28,2).Value = "rover finds evidence american moon"
$ws.Cells.Item(28,3).Value = 1
$ws.Cells.Item(29,1).Value = 5
$ws.Cells.Item(29,2).Value = "lunar rover finds evidence american"
$ws.Cells.Item(29,3).Value = 1
$ws.Cells.Item(30,1).Value = 5
$ws.Cells.Item(30,2).Value = "finds evidence american moon landings"
$ws.Cells.Item(30,3).Value = 1
$ws.Cells.Item(31,1).Value = 5
$ws.Cells.Item(31,2).Value = "chinese lunar rover finds evidence"
$ws.Cells.Item(31,3).Value = 1
Write-Host "ngrams updated"
